$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6346
$ws.Range("J3").Value = 6723
$ws.Range("J4").Value = 1457
$ws.Range("J6").Value = 8837
$ws.Range("J7").Value = 23878

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J6").Value = 178
$ws.Range("J7").Value = 695
$ws.Range("J8").Value = 1507
$ws.Range("J16").Value = 98
$ws.Range("J18").Value = 206
$ws.Range("J19").Value = 705
$ws.Range("J20").Value = 496
$ws.Range("J22").Value = 57
$ws.Range("J23").Value = 223
$ws.Range("J29").Value = 1310
$ws.Range("J31").Value = 220
$ws.Range("J33").Value = 1081
$ws.Range("J34").Value = 111
$ws.Range("J35").Value = 31
$ws.Range("J36").Value = 323
$ws.Range("J37").Value = 734
$ws.Range("J42").Value = 1031
$ws.Range("J45").Value = 35
$ws.Range("J47").Value = 179
$ws.Range("J48").Value = 277
$ws.Range("J49").Value = 152
$ws.Range("J51").Value = 297
$ws.Range("J53").Value = 345
$ws.Range("J54").Value = 457
$ws.Range("J55").Value = 346
$ws.Range("J56").Value = 33
$ws.Range("J59").Value = 28
$ws.Range("J63").Value = 85
$ws.Range("J65").Value = 588
$ws.Range("J67").Value = 900
$ws.Range("J68").Value = 52
$ws.Range("J73").Value = 231
$ws.Range("J76").Value = 359
$ws.Range("J78").Value = 283
$ws.Range("J79").Value = 676
$ws.Range("J83").Value = 477
$ws.Range("J85").Value = 995
$ws.Range("J86").Value = 153
$ws.Range("J89").Value = 311
$ws.Range("J96").Value = 263
$ws.Range("J101").Value = 23878

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 79
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 222
$ws.Range("J7").Value = 695

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 97
$ws.Range("J7").Value = 311

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 265
$ws.Range("J3").Value = 353
$ws.Range("J7").Value = 995

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 61
$ws.Range("J3").Value = 43
$ws.Range("J7").Value = 345

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 407
$ws.Range("J3").Value = 454
$ws.Range("J6").Value = 529
$ws.Range("J7").Value = 1507

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J4").Value = 18
$ws.Range("J7").Value = 477

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 357
$ws.Range("J4").Value = 45
$ws.Range("J7").Value = 1081

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 216
$ws.Range("J3").Value = 249
$ws.Range("J4").Value = 25
$ws.Range("J6").Value = 216
$ws.Range("J7").Value = 734

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 168
$ws.Range("J7").Value = 588

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 226
$ws.Range("J3").Value = 337
$ws.Range("J6").Value = 246
$ws.Range("J7").Value = 900

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J4").Value = 10
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 111
$ws.Range("J3").Value = 94
$ws.Range("J6").Value = 214
$ws.Range("J7").Value = 457

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 459
$ws.Range("J7").Value = 1310

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 47
$ws.Range("J3").Value = 51
$ws.Range("J7").Value = 277

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 167
$ws.Range("J6").Value = 273
$ws.Range("J7").Value = 705

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 198
$ws.Range("J7").Value = 359

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 178

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 552
$ws.Range("J7").Value = 1031

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 84
$ws.Range("J7").Value = 283

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 73
$ws.Range("J6").Value = 184
$ws.Range("J7").Value = 346

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J4").Value = 22
$ws.Range("J7").Value = 223

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 188
$ws.Range("J4").Value = 42
$ws.Range("J6").Value = 200
$ws.Range("J7").Value = 676

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 144
$ws.Range("J6").Value = 131
$ws.Range("J7").Value = 496

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 55
$ws.Range("J7").Value = 206

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 105
$ws.Range("J6").Value = 98
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J6").Value = 80
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 83
$ws.Range("J7").Value = 153

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 79
$ws.Range("J7").Value = 297

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("J2").Value = 21
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("J2").Value = 27
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J5").Value = 1
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 98
